# Adds a new forecast-vintage column (AF, dated 2020-05-11) and a new
# observation row (44, dated 2020-05-25) to both the "cases" and "deaths"
# sheets, mirroring the existing diagonal forecast-table layout.

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# Numeric values to drop into the new diagonal (row 31..43, column AF).
$afValuesBySheet = @{
    "cases"  = @{ 31 = 47802; 32 = 49139; 33 = 50314; 34 = 51294; 35 = 52445;
                  36 = 53505; 37 = 54586; 38 = 55608; 39 = 56587; 40 = 57758;
                  41 = 58548; 42 = 59280; 43 = 60046 }
    "deaths" = @{ 31 = 3887;  32 = 3997;  33 = 4101;  34 = 4198;  35 = 4297;
                  36 = 4389;  37 = 4468;  38 = 4549;  39 = 4630;  40 = 4716;
                  41 = 4798;  42 = 4874;  43 = 4951 }
}

# B30 was previously blank; the new vintage fills in the first diagonal cell.
$b30BySheet = @{ "cases" = 46131; "deaths" = 3743 }

# New row 44 (2020-05-25) only has one populated cell: AF44.
$af44BySheet = @{ "cases" = 60924; "deaths" = 5025 }

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- New column AF (col 32), header = "2020-05-11" ------------------
    # Force text so Excel's date auto-detection doesn't turn the literal
    # "2020-05-11" string into a date serial number.
    $headerCell = $ws.Range("AF1")
    $headerCell.NumberFormat = "@"
    $headerCell.Value = "2020-05-11"
    $headerCell.Style = "Normal"

    # Rows 2..29: column AF stays empty, but the cell needs to materialize
    # (matches the pre-existing pattern of explicit empty <c> cells).
    for ($r = 2; $r -le 29; $r++) {
        $ws.Cells.Item($r, 32).Style = "Normal"
    }

    # Row 30: B30 gains its first value; AF30 stays empty (diagonal hasn't
    # reached this row/column pair yet).
    $ws.Cells.Item(30, 2).Value = $b30BySheet[$sheetName]
    $ws.Cells.Item(30, 32).Style = "Normal"

    # Rows 31..43: column AF gets the new forecast values (diagonal).
    $afValues = $afValuesBySheet[$sheetName]
    foreach ($r in $afValues.Keys) {
        $ws.Cells.Item($r, 32).Value = $afValues[$r]
    }

    # --- New row 44 (2020-05-25) -----------------------------------------
    $rowHeaderCell = $ws.Cells.Item(44, 1)
    $rowHeaderCell.NumberFormat = "@"
    $rowHeaderCell.Value = "2020-05-25"
    $rowHeaderCell.Style = "Normal"

    # Columns B..AE (2..31) stay empty for row 44.
    for ($c = 2; $c -le 31; $c++) {
        $ws.Cells.Item(44, $c).Style = "Normal"
    }

    # AF44 is the only populated cell in the new row.
    $ws.Cells.Item(44, 32).Value = $af44BySheet[$sheetName]
}
